$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Delete the empty "Sheet1" worksheet.
$wb.Worksheets.Item("Sheet1").Delete()

# Rename the remaining data sheet ("18") to "13".
$ws = $wb.Worksheets.Item("18")
$ws.Name = "13"

# The chart's series formulas still reference the old sheet name -
# repoint them at the renamed sheet.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('13'!`$A`$12,'13'!`$B`$11:`$Q`$11,'13'!`$B`$12:`$Q`$12,1)"

# Update the active selection on the remaining sheet.
$ws.Activate()
$ws.Range("H15").Select()
